$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 4555
$ws.Range("J29").Value = 10821
$ws.Range("L29").Value = 32463
$ws.Range("N29").Value = -33025
$ws.Range("H47").Value = 16344.5
$ws.Range("I47").Value = 20016.75
$ws.Range("K47").Value = 20016.75
$ws.Range("M47").Value = -19044.75
$ws.Range("H64").Value = 16670833
$ws.Range("I64").Value = 33336666
$ws.Range("K64").Value = 33336666
$ws.Range("M64").Value = -33336418
$ws.Range("H67").Value = 16670833
$ws.Range("I67").Value = 33336666
$ws.Range("K67").Value = 33336666
$ws.Range("M67").Value = -33335808
$ws.Range("H69").Value = 8010
$ws.Range("J69").Value = 9015
$ws.Range("L69").Value = 27045
$ws.Range("N69").Value = -28793
$ws.Range("H72").Value = 8010
$ws.Range("J72").Value = 9015
$ws.Range("L72").Value = 81135
$ws.Range("N72").Value = -89871
$ws.Range("H76").Value = 3499.5
$ws.Range("I76").Value = 3499.5
$ws.Range("K76").Value = 3499.5
$ws.Range("M76").Value = -3184.5
$ws.Range("H79").Value = 3499.5
$ws.Range("I79").Value = 3499.5
$ws.Range("K79").Value = 3499.5
$ws.Range("M79").Value = -2407.5
$ws.Range("H132").Value = 2013.0769
$ws.Range("I132").Value = 2013.0769
$ws.Range("K132").Value = 6039.2307
$ws.Range("M132").Value = -3509.2307
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 49250
$ws.Range("J24").Value = 49250
$ws.Range("L24").Value = 49250
$ws.Range("N24").Value = -49998
$ws.Range("H100").Value = 49250
$ws.Range("J100").Value = 49250
$ws.Range("L100").Value = 49250
$ws.Range("N100").Value = -51414
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 317.875
$ws.Range("I80").Value = 319
$ws.Range("J80").Value = 316
$ws.Range("K80").Value = 319
$ws.Range("L80").Value = 316
$ws.Range("M80").Value = 679
$ws.Range("N80").Value = -2312
$ws.Range("H82").Value = 19252.334
$ws.Range("I82").Value = 19252.334
$ws.Range("K82").Value = 19252.334
$ws.Range("M82").Value = -18869.334
$ws.Range("H83").Value = 317.875
$ws.Range("I83").Value = 319
$ws.Range("J83").Value = 316
$ws.Range("K83").Value = 1595
$ws.Range("L83").Value = 1580
$ws.Range("M83").Value = 3397
$ws.Range("N83").Value = -11564
$ws.Range("H85").Value = 19252.334
$ws.Range("I85").Value = 19252.334
$ws.Range("K85").Value = 19252.334
$ws.Range("M85").Value = -17926.334
$ws.Range("H94").Value = 1292.625
$ws.Range("I94").Value = 1292.625
$ws.Range("K94").Value = 1292.625
$ws.Range("M94").Value = -841.625
$ws.Range("H95").Value = 39000
$ws.Range("J95").Value = 39000
$ws.Range("L95").Value = 39000
$ws.Range("N95").Value = -44492
$ws.Range("H97").Value = 14666.333
$ws.Range("I97").Value = 14264
$ws.Range("J97").Value = 15471
$ws.Range("K97").Value = 14264
$ws.Range("L97").Value = 15471
$ws.Range("M97").Value = -13273
$ws.Range("N97").Value = -17453
$ws.Range("H100").Value = 13569.25
$ws.Range("J100").Value = 13569.25
$ws.Range("L100").Value = 13569.25
$ws.Range("N100").Value = -15733.25
$ws.Range("H105").Value = 6500
$ws.Range("I105").Value = 6500
$ws.Range("K105").Value = 6500
$ws.Range("M105").Value = -4753
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H29").Value = 2019
$ws.Range("I29").Value = 2019
$ws.Range("K29").Value = 2019
$ws.Range("M29").Value = -1726
$ws.Range("H31").Value = 1690
$ws.Range("I31").Value = 1036.4
$ws.Range("K31").Value = 1036.4
$ws.Range("M31").Value = -741.4000000000001
$ws.Range("H34").Value = 1690
$ws.Range("I34").Value = 1036.4
$ws.Range("K34").Value = 1036.4
$ws.Range("M34").Value = -834.4000000000001
$ws.Range("H94").Value = 970.7143
$ws.Range("I94").Value = 698.75
$ws.Range("K94").Value = 698.75
$ws.Range("M94").Value = -247.75
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 273149.3
$ws.Range("I4").Value = 278229.56
$ws.Range("J4").Value = 250288.25
$ws.Range("K4").Value = 834688.6799999999
$ws.Range("L4").Value = 750864.75
$ws.Range("M4").Value = -834576.6799999999
$ws.Range("N4").Value = -751088.75
$ws.Range("H5").Value = 319.66666
$ws.Range("I5").Value = 345.875
$ws.Range("J5").Value = 110
$ws.Range("K5").Value = 1037.625
$ws.Range("L5").Value = 330
$ws.Range("M5").Value = -925.625
$ws.Range("N5").Value = -554
$ws.Range("H18").Value = 708.3333
$ws.Range("I18").Value = 708.3333
$ws.Range("K18").Value = 2124.9999
$ws.Range("M18").Value = -1955.9999
$ws.Range("H86").Value = 9567.166999999999
$ws.Range("I86").Value = 1100
$ws.Range("J86").Value = 26501.5
$ws.Range("K86").Value = 3300
$ws.Range("L86").Value = 79504.5
$ws.Range("M86").Value = -2114
$ws.Range("N86").Value = -81876.5
$ws.Range("H89").Value = 9567.166999999999
$ws.Range("I89").Value = 1100
$ws.Range("J89").Value = 26501.5
$ws.Range("K89").Value = 9900
$ws.Range("L89").Value = 238513.5
$ws.Range("M89").Value = -3972
$ws.Range("N89").Value = -250369.5
$ws.Range("H113").Value = 596
$ws.Range("J113").Value = 756
$ws.Range("L113").Value = 2268
$ws.Range("N113").Value = -6608
$ws.Range("H122").Value = 1219.6
$ws.Range("I122").Value = 700
$ws.Range("J122").Value = 1349.5
$ws.Range("K122").Value = 6300
$ws.Range("L122").Value = 12145.5
$ws.Range("M122").Value = -3850
$ws.Range("N122").Value = -17045.5
$ws.Range("H135").Value = 319.66666
$ws.Range("I135").Value = 345.875
$ws.Range("J135").Value = 110
$ws.Range("K135").Value = 3112.875
$ws.Range("L135").Value = 990
$ws.Range("M135").Value = -577.875
$ws.Range("N135").Value = -6060
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1000
$ws.Range("I102").Value = 1000
$ws.Range("K102").Value = 1000
$ws.Range("M102").Value = 622
$ws.Range("H130").Value = 99995
$ws.Range("J130").Value = 99995
$ws.Range("L130").Value = 99995
$ws.Range("N130").Value = -110035
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4210
$ws.Range("I46").Value = 2100
$ws.Range("J46").Value = 5001.25
$ws.Range("K46").Value = 2100
$ws.Range("L46").Value = 5001.25
$ws.Range("M46").Value = -1912
$ws.Range("N46").Value = -5377.25
$ws.Range("H82").Value = 2271.375
$ws.Range("I82").Value = 1940.3636
$ws.Range("K82").Value = 1940.3636
$ws.Range("M82").Value = -1579.3636
$ws.Range("H85").Value = 2271.375
$ws.Range("I85").Value = 1940.3636
$ws.Range("K85").Value = 1940.3636
$ws.Range("M85").Value = -692.3635999999999
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H132").Value = 7896.8
$ws.Range("J132").Value = 3992
$ws.Range("L132").Value = 11976
$ws.Range("N132").Value = -17036
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 20000
$ws.Range("J47").Value = 20000
$ws.Range("L47").Value = 20000
$ws.Range("N47").Value = -21144
$ws.Range("H96").Value = 5429.6
$ws.Range("I96").Value = 4787
$ws.Range("K96").Value = 4787
$ws.Range("M96").Value = -3414
$ws.Range("H122").Value = 7600.4
$ws.Range("I122").Value = 7749.25
$ws.Range("J122").Value = 7005
$ws.Range("K122").Value = 23247.75
$ws.Range("L122").Value = 21015
$ws.Range("M122").Value = -20797.75
$ws.Range("N122").Value = -25915
$ws.Range("H126").Value = 3766.5
$ws.Range("I126").Value = 4220.8
$ws.Range("K126").Value = 12662.4
$ws.Range("M126").Value = -10192.4
$ws.Range("H132").Value = 2666.3333
$ws.Range("I132").Value = 2859.6
$ws.Range("K132").Value = 8578.799999999999
$ws.Range("M132").Value = -6048.799999999999
